# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" 1804 rows need to move ahead of the 1805 rows for each
# worker, keeping each worker's Valor Mora / Salario Basico together with
# their period. Net effect: rows 16-19 get reordered to
#   16: ADOLFO JIMENEZ DIMAS   / 1804 / 1042  / 781242
#   17: CARLOS ALBERTO CERVANTES JULIO / 1804 / 1042  / 781242
#   18: ADOLFO JIMENEZ DIMAS   / 1805 / 31249 / 781242
#   19: CARLOS ALBERTO CERVANTES JULIO / 1805 / 31249 / 781242

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1804"
$ws.Range("F16").Value = 1042

$ws.Range("C17").Value = "1047377965"
$ws.Range("D17").Value = "CARLOS ALBERTO CERVANTES JULIO"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 1042

$ws.Range("C18").Value = "73209005"
$ws.Range("D18").Value = "ADOLFO JIMENEZ DIMAS"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 31249

$ws.Range("C19").Value = "1047377965"
$ws.Range("D19").Value = "CARLOS ALBERTO CERVANTES JULIO"
$ws.Range("E19").Value = "1805"
$ws.Range("F19").Value = 31249
